# Daily attendance processing - 2025-11-17 08:55:16
# Rotates the "Recorded By" (column G) list so that a trailing "System"
# entry (or, if absent, whatever entry is last) is moved to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "
        $n = $parts.Length

        if ($n -ge 2 -and -not $parts[0].Equals("System")) {
            $newParts = @($parts[$n - 1]) + $parts[0..($n - 2)]
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value = $newVal
        }
    }
}
